$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row above row 20 (shifting existing rows 20+ down by one)
$ws.Rows(20).Insert()

$ws.Range("B20").Value = "writeLong(slaveId,register,value) or writeLong([slaveId,register,value],..,[slaveId,register,value])"
$ws.Range("C20").Value = "write 32bit integer to two 16bit int registers: MODBUS function 16"
$ws.Rows(20).RowHeight = 13.8
